# Apply odds updates to the Jogos do Dia Betfair Back/Lay workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.9
$ws.Range("G2").Value = 4.6
$ws.Range("H2").Value = 1.81
$ws.Range("I2").Value = 1.98
$ws.Range("J2").Value = 3.95
$ws.Range("K2").Value = 4.5

# Row 3
$ws.Range("S3").Value = 3.05

# Row 4
$ws.Range("G4").Value = 3.3
$ws.Range("H4").Value = 2.3
$ws.Range("R4").Value = 1.62
$ws.Range("S4").Value = 2.52
$ws.Range("U4").Value = 2.72

# Row 5
$ws.Range("F5").Value = 2.04
$ws.Range("H5").Value = 1.63
$ws.Range("I5").Value = 5.1
$ws.Range("J5").Value = 3.25

# Row 9
$ws.Range("F9").Value = 2.02
$ws.Range("P9").Value = 1.45
$ws.Range("Q9").Value = 2.6

# Row 10
$ws.Range("F10").Value = 2.8
$ws.Range("G10").Value = 3.05
$ws.Range("H10").Value = 3.1
$ws.Range("I10").Value = 3.55
$ws.Range("J10").Value = 2.74
$ws.Range("K10").Value = 2.98
$ws.Range("P10").Value = 1.37
$ws.Range("Q10").Value = 3.05

# Row 11
$ws.Range("F11").Value = 1.93
$ws.Range("G11").Value = 2.24
$ws.Range("H11").Value = 4.4
$ws.Range("I11").Value = 5.9
$ws.Range("J11").Value = 3.1
$ws.Range("K11").Value = 3.65
$ws.Range("P11").Value = 1.53
